$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.033348942174757
$ws.Cells.Item(2, 4).Value = 1.042280535274918
$ws.Cells.Item(2, 5).Value = 1.050627189925885
$ws.Cells.Item(2, 6).Value = 1.055509705240836
$ws.Cells.Item(2, 9).Value = 1.034158759020832
$ws.Cells.Item(2, 10).Value = 1.038474157959292
$ws.Cells.Item(2, 11).Value = 1.045057667961704
$ws.Cells.Item(2, 12).Value = 1.053380942084291
$ws.Cells.Item(2, 13).Value = 1.058249966318539
$ws.Cells.Item(2, 14).Value = 1.039948910378246

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.034330382771139
$ws.Cells.Item(3, 4).Value = 1.043049099221802
$ws.Cells.Item(3, 5).Value = 1.051599461722723
$ws.Cells.Item(3, 6).Value = 1.056477108623294
$ws.Cells.Item(3, 9).Value = 1.034311384687696
$ws.Cells.Item(3, 10).Value = 1.039098112899147
$ws.Cells.Item(3, 11).Value = 1.045637142082738
$ws.Cells.Item(3, 12).Value = 1.054165286458485
$ws.Cells.Item(3, 13).Value = 1.059030432529311
$ws.Cells.Item(3, 14).Value = 1.040573751405683

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.034965817471187
$ws.Cells.Item(4, 4).Value = 1.043546569137339
$ws.Cells.Item(4, 5).Value = 1.05222970333967
$ws.Cells.Item(4, 6).Value = 1.057103907580313
$ws.Cells.Item(4, 9).Value = 1.034408824706331
$ws.Cells.Item(4, 10).Value = 1.039501634388372
$ws.Cells.Item(4, 11).Value = 1.046011597954687
$ws.Cells.Item(4, 12).Value = 1.054673305958044
$ws.Cells.Item(4, 13).Value = 1.059535648620599
$ws.Cells.Item(4, 14).Value = 1.040977845941707

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.035233043797169
$ws.Cells.Item(5, 4).Value = 1.043755741998521
$ws.Cells.Item(5, 5).Value = 1.052494923098805
$ws.Cells.Item(5, 6).Value = 1.057367610003806
$ws.Cells.Item(5, 9).Value = 1.03444947223664
$ws.Cells.Item(5, 10).Value = 1.039671221576712
$ws.Cells.Item(5, 11).Value = 1.046168897943363
$ws.Cells.Item(5, 12).Value = 1.054886995156008
$ws.Cells.Item(5, 13).Value = 1.059748088848261
$ws.Cells.Item(5, 14).Value = 1.041147673963307

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.035277917506446
$ws.Cells.Item(6, 4).Value = 1.043790865147357
$ws.Cells.Item(6, 5).Value = 1.052539470238581
$ws.Cells.Item(6, 6).Value = 1.0574118982683
$ws.Cells.Item(6, 9).Value = 1.034456278576532
$ws.Cells.Item(6, 10).Value = 1.03969969286634
$ws.Cells.Item(6, 11).Value = 1.046195302165895
$ws.Cells.Item(6, 12).Value = 1.054922881391439
$ws.Cells.Item(6, 13).Value = 1.059783761245637
$ws.Cells.Item(6, 14).Value = 1.041176185685431

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.034969387811292
$ws.Cells.Item(7, 4).Value = 1.043549363973444
$ws.Cells.Item(7, 5).Value = 1.052233246174616
$ws.Cells.Item(7, 6).Value = 1.057107430415791
$ws.Cells.Item(7, 9).Value = 1.034409369083884
$ws.Cells.Item(7, 10).Value = 1.039503900630901
$ws.Cells.Item(7, 11).Value = 1.04601370028247
$ws.Cells.Item(7, 12).Value = 1.054676160821944
$ws.Cells.Item(7, 13).Value = 1.059538487072439
$ws.Cells.Item(7, 14).Value = 1.040980115402561

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.03368054614868
$ws.Cells.Item(8, 4).Value = 1.042540241675745
$ws.Cells.Item(8, 5).Value = 1.050955541797798
$ws.Cells.Item(8, 6).Value = 1.055836472353252
$ws.Cells.Item(8, 9).Value = 1.034210612394826
$ws.Cells.Item(8, 10).Value = 1.038685071334617
$ws.Cells.Item(8, 11).Value = 1.045253607469034
$ws.Cells.Item(8, 12).Value = 1.053645911585388
$ws.Cells.Item(8, 13).Value = 1.058513686022883
$ws.Cells.Item(8, 14).Value = 1.040160123274755

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.031412359731079
$ws.Cells.Item(9, 4).Value = 1.040763300994947
$ws.Cells.Item(9, 5).Value = 1.048712668016438
$ws.Cells.Item(9, 6).Value = 1.053603247467022
$ws.Cells.Item(9, 9).Value = 1.033850297140171
$ws.Cells.Item(9, 10).Value = 1.037240546735714
$ws.Cells.Item(9, 11).Value = 1.043910418167545
$ws.Cells.Item(9, 12).Value = 1.051834325589952
$ws.Cells.Item(9, 13).Value = 1.056709450598644
$ws.Cells.Item(9, 14).Value = 1.038713547285232

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.02990223528264
$ws.Cells.Item(10, 4).Value = 1.039579594510025
$ws.Cells.Item(10, 5).Value = 1.047223272592382
$ws.Cells.Item(10, 6).Value = 1.052118780524098
$ws.Cells.Item(10, 9).Value = 1.033603337745488
$ws.Cells.Item(10, 10).Value = 1.036276472168579
$ws.Cells.Item(10, 11).Value = 1.043012452178951
$ws.Cells.Item(10, 12).Value = 1.050629245227095
$ws.Cells.Item(10, 13).Value = 1.055507760143365
$ws.Cells.Item(10, 14).Value = 1.037748103621631

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.029248816567518
$ws.Cells.Item(11, 4).Value = 1.039067269842546
$ws.Cells.Item(11, 5).Value = 1.046579750483607
$ws.Cells.Item(11, 6).Value = 1.051477035971391
$ws.Cells.Item(11, 9).Value = 1.033494806616595
$ws.Cells.Item(11, 10).Value = 1.03585877568433
$ws.Cells.Item(11, 11).Value = 1.042623039639365
$ws.Cells.Item(11, 12).Value = 1.050108071216006
$ws.Cells.Item(11, 13).Value = 1.054987696166408
$ws.Cells.Item(11, 14).Value = 1.037329813960464

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.029006179677952
$ws.Cells.Item(12, 4).Value = 1.038877005079126
$ws.Cells.Item(12, 5).Value = 1.046340928547306
$ws.Cells.Item(12, 6).Value = 1.05123882074816
$ws.Cells.Item(12, 9).Value = 1.033454253967133
$ws.Cells.Item(12, 10).Value = 1.035703588391589
$ws.Cells.Item(12, 11).Value = 1.042478307096105
$ws.Cells.Item(12, 12).Value = 1.049914580130366
$ws.Cells.Item(12, 13).Value = 1.05479456366693
$ws.Cells.Item(12, 14).Value = 1.037174406283967

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.029058222851549
$ws.Cells.Item(13, 4).Value = 1.038917815894762
$ws.Cells.Item(13, 5).Value = 1.046392147105199
$ws.Cells.Item(13, 6).Value = 1.051289911587134
$ws.Cells.Item(13, 9).Value = 1.033462963472425
$ws.Cells.Item(13, 10).Value = 1.035736878220266
$ws.Cells.Item(13, 11).Value = 1.042509356658689
$ws.Cells.Item(13, 12).Value = 1.049956080259148
$ws.Cells.Item(13, 13).Value = 1.054835989308705
$ws.Cells.Item(13, 14).Value = 1.03720774338802

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.029228758656408
$ws.Cells.Item(14, 4).Value = 1.039051541768351
$ws.Cells.Item(14, 5).Value = 1.046560005080636
$ws.Cells.Item(14, 6).Value = 1.051457341811862
$ws.Cells.Item(14, 9).Value = 1.033491459402193
$ws.Cells.Item(14, 10).Value = 1.035845948589058
$ws.Cells.Item(14, 11).Value = 1.042611077779926
$ws.Cells.Item(14, 12).Value = 1.050092075215713
$ws.Cells.Item(14, 13).Value = 1.054971730902871
$ws.Cells.Item(14, 14).Value = 1.037316968649246

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.029333840939539
$ws.Cells.Item(15, 4).Value = 1.039133939422217
$ws.Cells.Item(15, 5).Value = 1.046663455889094
$ws.Cells.Item(15, 6).Value = 1.051560521973367
$ws.Cells.Item(15, 9).Value = 1.033508984981138
$ws.Cells.Item(15, 10).Value = 1.035913145664838
$ws.Cells.Item(15, 11).Value = 1.042673739965809
$ws.Cells.Item(15, 12).Value = 1.050175878965325
$ws.Cells.Item(15, 13).Value = 1.055055371428649
$ws.Cells.Item(15, 14).Value = 1.037384261152581

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.029945610586553
$ws.Cells.Item(16, 4).Value = 1.039613600715329
$ws.Cells.Item(16, 5).Value = 1.047266010540331
$ws.Cells.Item(16, 6).Value = 1.052161392991221
$ws.Cells.Item(16, 9).Value = 1.033610507008779
$ws.Cells.Item(16, 10).Value = 1.036304188189988
$ws.Cells.Item(16, 11).Value = 1.043038283904972
$ws.Cells.Item(16, 12).Value = 1.050663847261696
$ws.Cells.Item(16, 13).Value = 1.055542280972249
$ws.Cells.Item(16, 14).Value = 1.03777585900297

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.030329485104183
$ws.Cells.Item(17, 4).Value = 1.039914541452087
$ws.Cells.Item(17, 5).Value = 1.047644351664625
$ws.Cells.Item(17, 6).Value = 1.052538582577127
$ws.Cells.Item(17, 9).Value = 1.033673762086652
$ws.Cells.Item(17, 10).Value = 1.03654941352106
$ws.Cells.Item(17, 11).Value = 1.043266796041657
$ws.Cells.Item(17, 12).Value = 1.050970107211516
$ws.Cells.Item(17, 13).Value = 1.055847781016315
$ws.Cells.Item(17, 14).Value = 1.03802143258214

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.030553438423583
$ws.Cells.Item(18, 4).Value = 1.040090097017764
$ws.Cells.Item(18, 5).Value = 1.047865166255065
$ws.Cells.Item(18, 6).Value = 1.052758691096052
$ws.Cells.Item(18, 9).Value = 1.033710503675289
$ws.Cells.Item(18, 10).Value = 1.036692425540383
$ws.Cells.Item(18, 11).Value = 1.043400026577678
$ws.Cells.Item(18, 12).Value = 1.051148804542462
$ws.Cells.Item(18, 13).Value = 1.056026000580938
$ws.Cells.Item(18, 14).Value = 1.038164647694934

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.030629808471763
$ws.Cells.Item(19, 4).Value = 1.040149960612318
$ws.Cells.Item(19, 5).Value = 1.04794048112331
$ws.Cells.Item(19, 6).Value = 1.05283375937721
$ws.Cells.Item(19, 9).Value = 1.033723005468402
$ws.Cells.Item(19, 10).Value = 1.036741184890166
$ws.Cells.Item(19, 11).Value = 1.043445445099597
$ws.Cells.Item(19, 12).Value = 1.051209746046795
$ws.Cells.Item(19, 13).Value = 1.056086773322207
$ws.Cells.Item(19, 14).Value = 1.038213476288586

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.030288294259367
$ws.Cells.Item(20, 4).Value = 1.039882251082791
$ws.Cells.Item(20, 5).Value = 1.047603745310844
$ws.Cells.Item(20, 6).Value = 1.052498103330445
$ws.Cells.Item(20, 9).Value = 1.033666991344113
$ws.Cells.Item(20, 10).Value = 1.036523105617009
$ws.Cells.Item(20, 11).Value = 1.043242284713365
$ws.Cells.Item(20, 12).Value = 1.050937242092352
$ws.Cells.Item(20, 13).Value = 1.055815000998839
$ws.Cells.Item(20, 14).Value = 1.037995087317847

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.029178538122884
$ws.Cells.Item(21, 4).Value = 1.039012161844214
$ws.Cells.Item(21, 5).Value = 1.046510569261867
$ws.Cells.Item(21, 6).Value = 1.051408033433181
$ws.Cells.Item(21, 9).Value = 1.033483074664993
$ws.Cells.Item(21, 10).Value = 1.035813831078084
$ws.Cells.Item(21, 11).Value = 1.042581125850514
$ws.Cells.Item(21, 12).Value = 1.05005202542893
$ws.Cells.Item(21, 13).Value = 1.054931757207984
$ws.Cells.Item(21, 14).Value = 1.037284805527722

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.028481206736565
$ws.Cells.Item(22, 4).Value = 1.038465307766889
$ws.Cells.Item(22, 5).Value = 1.045824466424179
$ws.Cells.Item(22, 6).Value = 1.050723573735511
$ws.Cells.Item(22, 9).Value = 1.033366054327454
$ws.Cells.Item(22, 10).Value = 1.035367672817713
$ws.Cells.Item(22, 11).Value = 1.042164923858511
$ws.Cells.Item(22, 12).Value = 1.049496011321965
$ws.Cells.Item(22, 13).Value = 1.054376672673261
$ws.Cells.Item(22, 14).Value = 1.036838013671446

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.028850835655364
$ws.Cells.Item(23, 4).Value = 1.038755185585723
$ws.Cells.Item(23, 5).Value = 1.046188066477351
$ws.Cells.Item(23, 6).Value = 1.051086332092906
$ws.Cells.Item(23, 9).Value = 1.033428220122298
$ws.Cells.Item(23, 10).Value = 1.035604209344945
$ws.Cells.Item(23, 11).Value = 1.042385608008583
$ws.Cells.Item(23, 12).Value = 1.049790711834502
$ws.Cells.Item(23, 13).Value = 1.054670909857692
$ws.Cells.Item(23, 14).Value = 1.037074886107678

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.030306906487534
$ws.Cells.Item(24, 4).Value = 1.039896841643101
$ws.Cells.Item(24, 5).Value = 1.047622093158388
$ws.Cells.Item(24, 6).Value = 1.052516393850042
$ws.Cells.Item(24, 9).Value = 1.033670051227173
$ws.Cells.Item(24, 10).Value = 1.036534993099908
$ws.Cells.Item(24, 11).Value = 1.043253360502872
$ws.Cells.Item(24, 12).Value = 1.050952092236367
$ws.Cells.Item(24, 13).Value = 1.055829812795741
$ws.Cells.Item(24, 14).Value = 1.038006991682335

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.031998390633717
$ws.Cells.Item(25, 4).Value = 1.041222525777532
$ws.Cells.Item(25, 5).Value = 1.049291477909878
$ws.Cells.Item(25, 6).Value = 1.054179828090162
$ws.Cells.Item(25, 9).Value = 1.033944638948899
$ws.Cells.Item(25, 10).Value = 1.037614180874963
$ws.Cells.Item(25, 11).Value = 1.044258110552632
$ws.Cells.Item(25, 12).Value = 1.052302202332692
$ws.Cells.Item(25, 13).Value = 1.057175693432831
$ws.Cells.Item(25, 14).Value = 1.039087712027816
